$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 and row 4 get refreshed "Correspond Handoff Datetime" (E)
# and "Correspond Handback DateTime" (H) timestamps.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-20 10:16:39"
$wsZh.Range("E4").Value = "2016-03-20 10:16:39"
$wsZh.Range("H3").Value = "2016-03-20 10:17:01"
$wsZh.Range("H4").Value = "2016-03-20 10:17:01"

# de-de sheet: same update, with its own timestamps.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-20 10:16:42"
$wsDe.Range("E4").Value = "2016-03-20 10:16:42"
$wsDe.Range("H3").Value = "2016-03-20 10:17:07"
$wsDe.Range("H4").Value = "2016-03-20 10:17:07"
